# Apply latest crypto market snapshot to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates as scraped: reference -> new value
$updates = [ordered]@{
    'D2' = '62.121.92'
    'E2' = '  +2.49%  '
    'D3' = '2.427.31'
    'E3' = '  +4.41%  '
    'E4' = '  -1.14%  '
    'D5' = '556.92'
    'E5' = '  +2.89%  '
    'D6' = '144.08'
    'E6' = '  +6.57%  '
    'D7' = '0.998'
    'E7' = '  +0.54%  '
    'D8' = '0.531'
    'E8' = '  +1.68%  '
    'D9' = '2.433.43'
    'E9' = '  +4.04%  '
    'D10' = '0.109'
    'E10' = '  +5.05%  '
    'E11' = '  +1.39%  '
    'D12' = '5.40'
    'E12' = '  +2.42%  '
    'D13' = '0.355'
    'E13' = '  +4.71%  '
    'D14' = '26.31'
    'E14' = '  +8.00%  '
    'D15' = '0.0000176'
    'E15' = '  +10.81%  '
    'D16' = '2.863.41'
    'E16' = '  +3.90%  '
    'D17' = '61.876.12'
    'E17' = '  +25.16%  '
    'D18' = '2.431.92'
    'E18' = '  +9.26%  '
    'D19' = '11.20'
    'E19' = '  +6.48%  '
    'B20' = 'Polkadot'
    'C20' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D20' = '4.20'
    'E20' = '  +3.32%  '
    'B21' = 'BitcoinCash'
    'C21' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D21' = '324.62'
    'E21' = '  +3.58%  '
    'D22' = '6.76'
    'E22' = '  +3.85%  '
    'E23' = '  +0.60%  '
    'D24' = '64.79'
    'E24' = '  +3.56%  '
    'D25' = '1.76'
    'E25' = '  +4.12%  '
    'D26' = '9.16'
    'E26' = '  +10.24%  '
    'D27' = '566.52'
    'E27' = '  +13.67%  '
    'D28' = '0.997'
    'E28' = '  -0.08%  '
    'D29' = '2.520.35'
    'E29' = '  +2.86%  '
    'B30' = 'PEPE'
    'C30' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D30' = '0.0₃0943'
    'E30' = '  +9.14%  '
    'B31' = 'InternetComputer(DFINITY)'
    'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D31' = '8.36'
    'E31' = '  +6.66%  '
    'D32' = '1.46'
    'E32' = '  +6.81%  '
    'D33' = '0.149'
    'E33' = '  +3.27%  '
    'D34' = '1.86'
    'E34' = '  +5.03%  '
    'D35' = '1.57'
    'E35' = '  +4.68%  '
    'D36' = '5.86'
    'E36' = '  +13.42%  '
    'B37' = 'Stacks'
    'C37' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D37' = '1.97'
    'E37' = '  +12.81%  '
    'B38' = 'FirstDigitalUSD'
    'C38' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D38' = '0.997'
    'E38' = '  +0.14%  '
    'D39' = '4.81'
    'E39' = '  +5.41%  '
    'D40' = '0.386'
    'E40' = '  +4.21%  '
    'D41' = '18.76'
    'E41' = '  +1.91%  '
    'D42' = '146.48'
    'E42' = '  +4.22%  '
    'E43' = '  +0.02%  '
    'B44' = 'dogwifhat'
    'C44' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D44' = '2.29'
    'E44' = '  +11.49%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '149.78'
    'E45' = '  +7.61%  '
    'D46' = '3.65'
    'E46' = '  +3.75%  '
    'D47' = '0.0541'
    'E47' = '  +6.42%  '
    'D48' = '20.40'
    'E48' = '  +7.00%  '
    'D49' = '0.596'
    'E49' = '  +5.51%  '
    'D50' = '0.0915'
    'E50' = '  +2.71%  '
    'D51' = '0.0226'
    'E51' = '  +3.45%  '
}

# Price cells whose text would otherwise be re-interpreted as a number
# (and lose formatting such as trailing zeros); force them to remain text.
$textCells = @('D5', 'D6', 'D7', 'D8', 'D10', 'D12', 'D13', 'D14', 'D15', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
